# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used four emoji values (📘 📕 📙 📗) as status
# markers. Replace them with plain, non-emoji markers:
#   📘 -> ⚠️
#   📕 -> -3
#   📙 -> +3
#   📗 -> ✅
#
# "-3" and "+3" must remain textual values (just like the emoji they
# replace), not be auto-converted to numbers by Excel. Prefixing the
# replacement text with a leading apostrophe forces Excel to store it
# as text (the standard "quote prefix" mechanism), exactly as typing
# '-3 into a cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("📘", "⚠️")
$ws.Cells.Replace("📕", "'-3")
$ws.Cells.Replace("📙", "'+3")
$ws.Cells.Replace("📗", "✅")
